# Scheduled-runner market data refresh for the job-leve profit sheets.
# Columns H:N on every sheet are plain, formula-free values (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) that the runner
# repopulates from the latest market-board pull; this mirrors that refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3219.6
$ws.Range("I29").Value = 2333.5
$ws.Range("J29").Value = 4548.75
$ws.Range("K29").Value = 7000.5
$ws.Range("L29").Value = 13646.25
$ws.Range("M29").Value = -6719.5
$ws.Range("N29").Value = -14208.25
$ws.Range("H38").Value = 5050
$ws.Range("J38").Value = 6105.846
$ws.Range("L38").Value = 18317.538
$ws.Range("N38").Value = -19061.538
$ws.Range("H58").Value = 956.2
$ws.Range("I58").Value = 821.25
$ws.Range("J58").Value = 1496
$ws.Range("K58").Value = 2463.75
$ws.Range("L58").Value = 4488
$ws.Range("M58").Value = -2313.75
$ws.Range("N58").Value = -4788
$ws.Range("H98").Value = 1289.5714
$ws.Range("J98").Value = 1249.3334
$ws.Range("L98").Value = 1249.3334
$ws.Range("N98").Value = -4245.3334
$ws.Range("H122").Value = 1289.5714
$ws.Range("J122").Value = 1249.3334
$ws.Range("L122").Value = 3748.0002
$ws.Range("N122").Value = -8648.0002
$ws.Range("H137").Value = 8585133
$ws.Range("I137").Value = 834456.0600000001
$ws.Range("K137").Value = 2503368.18
$ws.Range("M137").Value = -2500818.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13115.731
$ws.Range("I2").Value = 16940.645
$ws.Range("K2").Value = 16940.645
$ws.Range("M2").Value = -16827.645
$ws.Range("H37").Value = 42000
$ws.Range("J37").Value = 42000
$ws.Range("L37").Value = 42000
$ws.Range("N37").Value = -42546
$ws.Range("H61").Value = 5292.5454
$ws.Range("I61").Value = 6177.222
$ws.Range("J61").Value = 4680.077
$ws.Range("K61").Value = 6177.222
$ws.Range("L61").Value = 4680.077
$ws.Range("M61").Value = -5965.222
$ws.Range("N61").Value = -5104.077
$ws.Range("H74").Value = 10418368
$ws.Range("I74").Value = 13890402
$ws.Range("J74").Value = 2266
$ws.Range("K74").Value = 13890402
$ws.Range("L74").Value = 2266
$ws.Range("M74").Value = -13889528
$ws.Range("N74").Value = -4014
$ws.Range("H77").Value = 10418368
$ws.Range("I77").Value = 13890402
$ws.Range("J77").Value = 2266
$ws.Range("K77").Value = 69452010
$ws.Range("L77").Value = 11330
$ws.Range("M77").Value = -69447642
$ws.Range("N77").Value = -20066
$ws.Range("H102").Value = 2335.9285
$ws.Range("I102").Value = 2035.3
$ws.Range("J102").Value = 3087.5
$ws.Range("K102").Value = 2035.3
$ws.Range("L102").Value = 3087.5
$ws.Range("M102").Value = -413.3
$ws.Range("N102").Value = -6331.5
$ws.Range("H116").Value = 13115.731
$ws.Range("I116").Value = 16940.645
$ws.Range("K116").Value = 16940.645
$ws.Range("M116").Value = -14646.645
$ws.Range("H136").Value = 5292.5454
$ws.Range("I136").Value = 6177.222
$ws.Range("J136").Value = 4680.077
$ws.Range("K136").Value = 18531.666
$ws.Range("L136").Value = 14040.231
$ws.Range("M136").Value = -15981.666
$ws.Range("N136").Value = -19140.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13115.731
$ws.Range("I3").Value = 16940.645
$ws.Range("K3").Value = 16940.645
$ws.Range("M3").Value = -16826.645
$ws.Range("H20").Value = 3207.1707
$ws.Range("J20").Value = 3058.9583
$ws.Range("L20").Value = 3058.9583
$ws.Range("N20").Value = -3552.9583
$ws.Range("H35").Value = 66200
$ws.Range("J35").Value = 66200
$ws.Range("L35").Value = 66200
$ws.Range("N35").Value = -66820
$ws.Range("H99").Value = 2473.6553
$ws.Range("I99").Value = 2617
$ws.Range("J99").Value = 2097.375
$ws.Range("K99").Value = 2617
$ws.Range("L99").Value = 2097.375
$ws.Range("M99").Value = -1119
$ws.Range("N99").Value = -5093.375
$ws.Range("H134").Value = 2332.0637
$ws.Range("I134").Value = 2382.2273
$ws.Range("J134").Value = 1596.3334
$ws.Range("K134").Value = 7146.6819
$ws.Range("L134").Value = 4789.0002
$ws.Range("M134").Value = -4611.6819
$ws.Range("N134").Value = -9859.0002
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14495590
$ws.Range("I31").Value = 21278534
$ws.Range("J31").Value = 4755.136
$ws.Range("K31").Value = 21278534
$ws.Range("L31").Value = 4755.136
$ws.Range("M31").Value = -21278239
$ws.Range("N31").Value = -5345.136
$ws.Range("H34").Value = 14495590
$ws.Range("I34").Value = 21278534
$ws.Range("J34").Value = 4755.136
$ws.Range("K34").Value = 21278534
$ws.Range("L34").Value = 4755.136
$ws.Range("M34").Value = -21278332
$ws.Range("N34").Value = -5159.136
$ws.Range("H58").Value = 2483.647
$ws.Range("I58").Value = 2226.913
$ws.Range("J58").Value = 3020.4546
$ws.Range("K58").Value = 2226.913
$ws.Range("L58").Value = 3020.4546
$ws.Range("M58").Value = -2023.913
$ws.Range("N58").Value = -3426.4546
$ws.Range("H107").Value = 516.65
$ws.Range("I107").Value = 484.33334
$ws.Range("K107").Value = 484.33334
$ws.Range("M107").Value = 1435.66666
$ws.Range("H122").Value = 3123.25
$ws.Range("I122").Value = 2998
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8994
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -6544
$ws.Range("N122").Value = -16900
$ws.Range("H136").Value = 2483.647
$ws.Range("I136").Value = 2226.913
$ws.Range("J136").Value = 3020.4546
$ws.Range("K136").Value = 6680.739
$ws.Range("L136").Value = 9061.363799999999
$ws.Range("M136").Value = -4130.739
$ws.Range("N136").Value = -14161.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1312.5
$ws.Range("J34").Value = 2491.4285
$ws.Range("L34").Value = 7474.2855
$ws.Range("N34").Value = -7642.2855
$ws.Range("H39").Value = 3324.4443
$ws.Range("J39").Value = 3353.5715
$ws.Range("L39").Value = 10060.7145
$ws.Range("N39").Value = -10648.7145
$ws.Range("H56").Value = 71436920
$ws.Range("I56").Value = 71436920
$ws.Range("K56").Value = 71436920
$ws.Range("M56").Value = -71436390
$ws.Range("H124").Value = 3247.6191
$ws.Range("I124").Value = 2440.4
$ws.Range("J124").Value = 3499.875
$ws.Range("K124").Value = 7321.200000000001
$ws.Range("L124").Value = 10499.625
$ws.Range("M124").Value = -2411.200000000001
$ws.Range("N124").Value = -20319.625
$ws.Range("H127").Value = 1183.375
$ws.Range("J127").Value = 1183.375
$ws.Range("L127").Value = 3550.125
$ws.Range("N127").Value = -13470.125
$ws.Range("H129").Value = 1707.375
$ws.Range("I129").Value = 1308.4286
$ws.Range("J129").Value = 4500
$ws.Range("K129").Value = 3925.2858
$ws.Range("L129").Value = 13500
$ws.Range("M129").Value = 1074.7142
$ws.Range("N129").Value = -23500
$ws.Range("H131").Value = 13597227
$ws.Range("J131").Value = 17625694
$ws.Range("L131").Value = 52877082
$ws.Range("N131").Value = -52887162
$ws.Range("H132").Value = 1857.5483
$ws.Range("J132").Value = 2448.5334
$ws.Range("L132").Value = 22036.8006
$ws.Range("N132").Value = -27096.8006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 45124.75
$ws.Range("J43").Value = 45124.75
$ws.Range("L43").Value = 45124.75
$ws.Range("N43").Value = -45426.75
$ws.Range("H46").Value = 20119.9
$ws.Range("J46").Value = 24999.875
$ws.Range("L46").Value = 24999.875
$ws.Range("N46").Value = -25311.875
$ws.Range("H57").Value = 20998.25
$ws.Range("J57").Value = 20998.25
$ws.Range("L57").Value = 20998.25
$ws.Range("N57").Value = -22638.25
$ws.Range("H102").Value = 12147567
$ws.Range("I102").Value = 15004850
$ws.Range("K102").Value = 15004850
$ws.Range("M102").Value = -15003228
$ws.Range("H132").Value = 145844
$ws.Range("I132").Value = 334616.16
$ws.Range("K132").Value = 1003848.48
$ws.Range("M132").Value = -1001318.48
$ws.Range("H136").Value = 66545.35000000001
$ws.Range("J136").Value = 66545.35000000001
$ws.Range("L136").Value = 199636.05
$ws.Range("N136").Value = -204736.05
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5002
$ws.Range("I7").Value = 3657.3845
$ws.Range("K7").Value = 3657.3845
$ws.Range("M7").Value = -3545.3845
$ws.Range("H40").Value = 15690807
$ws.Range("I40").Value = 12502866
$ws.Range("J40").Value = 18524532
$ws.Range("K40").Value = 12502866
$ws.Range("L40").Value = 18524532
$ws.Range("M40").Value = -12502730
$ws.Range("N40").Value = -18524804
$ws.Range("H55").Value = 620.1177
$ws.Range("J55").Value = 962.55554
$ws.Range("L55").Value = 962.55554
$ws.Range("N55").Value = -1308.55554
$ws.Range("H126").Value = 5002
$ws.Range("I126").Value = 3657.3845
$ws.Range("K126").Value = 10972.1535
$ws.Range("M126").Value = -8502.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1566.0834
$ws.Range("I100").Value = 1256.2858
$ws.Range("K100").Value = 2512.5716
$ws.Range("M100").Value = -1971.5716
$ws.Range("H107").Value = 843.82355
$ws.Range("I107").Value = 635.4
$ws.Range("K107").Value = 1906.2
$ws.Range("M107").Value = 13.80000000000018
$ws.Range("H132").Value = 1969.6727
$ws.Range("J132").Value = 3322.6
$ws.Range("L132").Value = 9967.799999999999
$ws.Range("N132").Value = -15027.8
